# Mifos -> Finflux: 1st changes
#
# The "Repayment schedule" sheet gains a new (empty) column inserted
# before the old "Late" column (old N), which pushes the old N/O/P
# columns (Late / <blank> / Outstanding) one slot to the right
# (-> O/P/Q). The sheet also becomes the active/selected sheet with a
# new selection, and the previously-active "Edit Repayment Schedule"
# sheet loses its tab-selected flag automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N,O,P -> O,P,Q).
$ws.Columns("N").Insert() | Out-Null

# Make "Repayment schedule" the active sheet/tab and move the
# selection to L15, matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("L15").Select() | Out-Null
